# Apply the SMTP Details / Email Details block that was added below the
# existing "Database Access Details" section, plus the minor workbook-view
# bookkeeping tweak that came along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the existing plain "s=1" formatting block (rows 29:32, the
#     Database Access Details section) down onto the new rows 34:41 so the
#     new cells pick up the same fill/border style used throughout the
#     sheet, before we put any values into them. ---
$ws.Range("A29:C32").Copy()
$ws.Range("A34:C37").PasteSpecial(-4122)
$ws.Range("A29:C32").Copy()
$ws.Range("A38:C41").PasteSpecial(-4122)

# --- SMTP Details section ---
$ws.Range("A34").Value = "SMTP Details"

$ws.Range("B35").Value = "Host"
$ws.Range("C35").Value = "smtp.hostinger.com"

$ws.Range("B36").Value = "username"
$ws.Range("C36").Value = "info@parastoneglobal-ksa.com"

$ws.Range("B37").Value = "password"
$ws.Range("C37").Value = "info@paraSKSA321#"

# Row 38 is left as a blank styled separator row (like the other section
# breaks), so nothing further is written there.

# --- Email Details section ---
$ws.Range("A39").Value = "Email Details"

$ws.Range("B40").Value = "email"
$ws.Range("C40").Value = "info@parastoneglobal-ksa.com"

$ws.Range("B41").Value = "password"
$ws.Range("C41").Value = "info@paraSKSA321#"

# --- Update the selection to match the author's saved cursor position ---
$ws.Range("C43").Select() | Out-Null
